$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Name = "Marks for Sub. Code 150"
$ws.Range("A1").Value = "Name"
$ws.Range("B1").Value = "Marks"
$ws.Range("C1").Value = "Grade"
$ws.Range("A2").Value = "Aarush"
